$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 59.4
$ws.Range("N2").Value = 54.77309453746771

$ws.Range("K3").Value = 53.8
$ws.Range("N3").Value = 54.77309453746771
